$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must stay text (matches source formatting,
# e.g. trailing zeros / specific decimal precision). Mark them as Text format first,
# via a single Union range (applied area-by-area) so the same style record is reused.
$textCells = $ws.Range("D5")
$textCells = $excel.Union($textCells, $ws.Range("D6"))
$textCells = $excel.Union($textCells, $ws.Range("D7"))
$textCells = $excel.Union($textCells, $ws.Range("D8"))
$textCells = $excel.Union($textCells, $ws.Range("D9"))
$textCells = $excel.Union($textCells, $ws.Range("D10"))
$textCells = $excel.Union($textCells, $ws.Range("D11"))
$textCells = $excel.Union($textCells, $ws.Range("D12"))
$textCells = $excel.Union($textCells, $ws.Range("D14"))
$textCells = $excel.Union($textCells, $ws.Range("D15"))
$textCells = $excel.Union($textCells, $ws.Range("D16"))
$textCells = $excel.Union($textCells, $ws.Range("D19"))
$textCells = $excel.Union($textCells, $ws.Range("D20"))
$textCells = $excel.Union($textCells, $ws.Range("D24"))
$textCells = $excel.Union($textCells, $ws.Range("D25"))
$textCells = $excel.Union($textCells, $ws.Range("D26"))
$textCells = $excel.Union($textCells, $ws.Range("D28"))
$textCells = $excel.Union($textCells, $ws.Range("D30"))
$textCells = $excel.Union($textCells, $ws.Range("D31"))
$textCells = $excel.Union($textCells, $ws.Range("D32"))
$textCells = $excel.Union($textCells, $ws.Range("D33"))
$textCells = $excel.Union($textCells, $ws.Range("D34"))
$textCells = $excel.Union($textCells, $ws.Range("D35"))
$textCells = $excel.Union($textCells, $ws.Range("D36"))
$textCells = $excel.Union($textCells, $ws.Range("D37"))
$textCells = $excel.Union($textCells, $ws.Range("D38"))
$textCells = $excel.Union($textCells, $ws.Range("D39"))
$textCells = $excel.Union($textCells, $ws.Range("D40"))
$textCells = $excel.Union($textCells, $ws.Range("D41"))
$textCells = $excel.Union($textCells, $ws.Range("D42"))
$textCells = $excel.Union($textCells, $ws.Range("D43"))
$textCells = $excel.Union($textCells, $ws.Range("D44"))
$textCells = $excel.Union($textCells, $ws.Range("D45"))
$textCells = $excel.Union($textCells, $ws.Range("D46"))
$textCells = $excel.Union($textCells, $ws.Range("D47"))
$textCells = $excel.Union($textCells, $ws.Range("D48"))
$textCells = $excel.Union($textCells, $ws.Range("D49"))
$textCells = $excel.Union($textCells, $ws.Range("D50"))
$textCells = $excel.Union($textCells, $ws.Range("D51"))
foreach ($area in $textCells.Areas) {
    $area.NumberFormat = "@"
}

# Apply the updated values from the crypto price refresh.
$ws.Range("D2").Value = "30.210.90"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.863.72"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "234.87"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.4674"
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").Value = "0.2830"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").Value = "0.06534"
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("D10").Value = "21.38"
$ws.Range("E10").Value = "  +2.72%  "
$ws.Range("D11").Value = "0.07853"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").Value = "97.46"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "1.868.37"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "5.101"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "0.6725"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "279.35"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").Value = "30.195.16"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "5.518"
$ws.Range("E19").Value = "  +2.08%  "
$ws.Range("D20").Value = "12.63"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").Value = "2.110.12"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "6.147"
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("D25").Value = "9.199"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("D26").Value = "164.66"
$ws.Range("E26").Value = "  -1.54%  "
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "1.924"
$ws.Range("E28").Value = "  -3.63%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "0.09694"
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("D31").Value = "4.417"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("D32").Value = "1.476"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").Value = "4.080"
$ws.Range("E33").Value = "  -1.68%  "
$ws.Range("D34").Value = "0.04687"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "1.114"
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("D36").Value = "0.7052"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").Value = "2.729"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").Value = "0.01854"
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("D39").Value = "2.532"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "6.226"
$ws.Range("E40").Value = "  -7.52%  "
$ws.Range("D41").Value = "73.32"
$ws.Range("E41").Value = "  +1.36%  "
$ws.Range("D42").Value = "1.937"
$ws.Range("E42").Value = "  -1.83%  "
$ws.Range("D43").Value = "0.8465"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("D44").Value = "103.97"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.4163"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "7.200"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").Value = "9.159"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("D49").Value = "935.99"
$ws.Range("E49").Value = "  -6.63%  "
$ws.Range("D50").Value = "34.04"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "0.1127"
$ws.Range("E51").Value = "  -2.17%  "
